# Apply cryptos list update (price/volume refresh + Avalanche/WrappedliquidstakedEther2.0 row swap)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = "29.339.16"
$ws.Range('D3').Value = "1.878.75"
$ws.Range('E3').Value = "  +0.37%  "
$ws.Range('E4').Value = "  +0.21%  "
$ws.Range('D5').Value = "'0.7253"
$ws.Range('E5').Value = "  +2.38%  "
$ws.Range('D6').Value = "'242.74"
$ws.Range('E6').Value = "  +0.44%  "
$ws.Range('E7').Value = "  +0.19%  "
$ws.Range('D8').Value = "'0.08009"
$ws.Range('E8').Value = "  +2.91%  "
$ws.Range('E9').Value = "  +2.29%  "
$ws.Range('D10').Value = "'25.00"
$ws.Range('E10').Value = "  -0.11%  "
$ws.Range('D11').Value = "'0.08226"
$ws.Range('E11').Value = "  -2.01%  "
$ws.Range('D12').Value = "1.885.56"
$ws.Range('E12').Value = "  +0.94%  "
$ws.Range('D13').Value = "'94.73"
$ws.Range('E13').Value = "  +4.06%  "
$ws.Range('D14').Value = "'5.230"
$ws.Range('E14').Value = "  -0.26%  "
$ws.Range('D15').Value = "'0.7126"
$ws.Range('E15').Value = "  +0.18%  "
$ws.Range('D16').Value = "'6.400"
$ws.Range('E16').Value = "  +5.42%  "
$ws.Range('D17').Value = "'0.000008508"
$ws.Range('E17').Value = "  +3.93%  "
$ws.Range('D18').Value = "29.335.44"
$ws.Range('E18').Value = "  +0.11%  "
$ws.Range('D19').Value = "'243.39"
$ws.Range('E19').Value = "  +1.47%  "
$ws.Range('B20').Value = "WrappedliquidstakedEther2.0"
$ws.Range('C20').Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range('D20').Value = "2.140.74"
$ws.Range('E20').Value = "  +1.14%  "
$ws.Range('B21').Value = "Avalanche"
$ws.Range('C21').Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range('D21').Value = "'13.25"
$ws.Range('E21').Value = "  +0.36%  "
$ws.Range('E22').Value = "  +0.22%  "
$ws.Range('D23').Value = "'7.774"
$ws.Range('E23').Value = "  +0.26%  "
$ws.Range('D24').Value = "'1.001"
$ws.Range('E24').Value = "  +0.06%  "
$ws.Range('D25').Value = "'0.1602"
$ws.Range('E25').Value = "  +0.78%  "
$ws.Range('D26').Value = "'162.60"
$ws.Range('E26').Value = "  -0.39%  "
$ws.Range('D27').Value = "'9.039"
$ws.Range('E27').Value = "  +0.22%  "
$ws.Range('D28').Value = "'18.54"
$ws.Range('E28').Value = "  +0.46%  "
$ws.Range('E29').Value = "  +0.01%  "
$ws.Range('D30').Value = "'4.412"
$ws.Range('E30').Value = "  +0.24%  "
$ws.Range('D31').Value = "'4.306"
$ws.Range('E31').Value = "  +0.10%  "
$ws.Range('D32').Value = "'1.191"
$ws.Range('E32').Value = "  -7.62%  "
$ws.Range('D33').Value = "'0.05362"
$ws.Range('E33').Value = "  +0.50%  "
$ws.Range('E34').Value = "  +0.17%  "
$ws.Range('D35').Value = "'0.7589"
$ws.Range('E35').Value = "  +1.87%  "
$ws.Range('D36').Value = "'1.177"
$ws.Range('E36').Value = "  +0.03%  "
$ws.Range('D37').Value = "'2.699"
$ws.Range('E37').Value = "  +0.05%  "
$ws.Range('D38').Value = "'0.01880"
$ws.Range('E38').Value = "  +0.54%  "
$ws.Range('D39').Value = "1.278.00"
$ws.Range('E39').Value = "  +3.85%  "
$ws.Range('D40').Value = "'2.760"
$ws.Range('E40').Value = "  +1.31%  "
$ws.Range('D41').Value = "'6.430"
$ws.Range('E41').Value = "  -2.05%  "
$ws.Range('D42').Value = "'113.27"
$ws.Range('E42').Value = "  +3.31%  "
$ws.Range('D43').Value = "'0.9067"
$ws.Range('E43').Value = "  +2.46%  "
$ws.Range('D44').Value = "'74.27"
$ws.Range('E44').Value = "  +2.57%  "
$ws.Range('E45').Value = "  +9.28%  "
$ws.Range('D46').Value = "'1.002"
$ws.Range('E46').Value = "  +0.21%  "
$ws.Range('D47').Value = "2.033.22"
$ws.Range('E47').Value = "  +0.88%  "
$ws.Range('D48').Value = "'0.5231"
$ws.Range('E48').Value = "  +0.73%  "
$ws.Range('D49').Value = "'1.796"
$ws.Range('E49').Value = "  +0.10%  "
$ws.Range('D50').Value = "'9.489"
$ws.Range('E50').Value = "  +0.96%  "
$ws.Range('D51').Value = "'0.4350"
$ws.Range('E51').Value = "  +0.89%  "
